# Update column A (取得日時 / timestamp) for rows 2-22 on the first sheet
# from "2025-12-08 18:28:28" to "2025-12-08 18:36:57".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025-12-08 18:28:28") {
        $cell.Value = "2025-12-08 18:36:57"
    }
}
